$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1 - copy header formatting (bold, border, alignment) from A1,
# then set its text.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "nome"

# New row 2 data. Values are stored as plain text (not numbers), so force a
# text number format before assigning numeric-looking strings, then clear the
# formatting again so the cell keeps the default (unstyled) appearance while
# remaining text.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "48998418335"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").ClearFormats()

$ws.Range("C2").Value = "andy "

# Row 3 update: B3 text changes from "1.0" to "1"; add C3 "sla"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1"
$ws.Range("B3").ClearFormats()

$ws.Range("C3").Value = "sla"
